$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N3").Value = 7.5
$ws.Range("W3").Value = 7.5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q7").Value = 1.72
$ws.Range("U7").Value = 1.63
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.87
$ws.Range("U9").Value = 1.77
$ws.Range("V9").Value = 1.87
$ws.Range("O10").Value = 1.13
$ws.Range("U10").Value = 2.28
$ws.Range("V10").Value = 1.6
$ws.Range("Q11").Value = 1.5
$ws.Range("O12").Value = 1.18
$ws.Range("P12").Value = 4.5
$ws.Range("Q12").Value = 1.65
$ws.Range("R12").Value = 2.2
$ws.Range("N14").Value = 29
$ws.Range("Q14").Value = 1.3
$ws.Range("R14").Value = 3.5
$ws.Range("G16").Value = 2.38
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("Q23").Value = 2.25
$ws.Range("R23").Value = 1.62
$ws.Range("G26").Value = 2.4
$ws.Range("I26").Value = 2.75
$ws.Range("L26").Value = 3.25
$ws.Range("N26").Value = 13
$ws.Range("AE26").Value = 13
$ws.Range("AH26").Value = 10
$ws.Range("AM26").Value = 26
$ws.Range("G27").Value = 1.95
$ws.Range("I27").Value = 4.2
$ws.Range("J27").Value = 2.75
$ws.Range("L27").Value = 4.75
$ws.Range("M27").Value = 1.1
$ws.Range("N27").Value = 7
$ws.Range("Z27").Value = 17
$ws.Range("AD27").Value = 6
$ws.Range("AH27").Value = 9.5
$ws.Range("AK27").Value = 41
$ws.Range("AX27").Value = 23
$ws.Range("G30").Value = 2
$ws.Range("I30").Value = 3.75
$ws.Range("Q30").Value = 1.98
$ws.Range("R30").Value = 1.88
$ws.Range("AA30").Value = 17
$ws.Range("AH30").Value = 11
$ws.Range("AQ30").Value = 41
$ws.Range("AX30").Value = 19
$ws.Range("G32").Value = 2.25
$ws.Range("L32").Value = 3.25
$ws.Range("AR32").Value = 41
$ws.Range("BD32").Value = 151
$ws.Range("G38").Value = 2.2
$ws.Range("H38").Value = 3.25
$ws.Range("I38").Value = 3.25
$ws.Range("J38").Value = 3
$ws.Range("X38").Value = 9.5
$ws.Range("Y38").Value = 9.5
$ws.Range("AD38").Value = 6.5
$ws.Range("AE38").Value = 19
$ws.Range("AN38").Value = 4
$ws.Range("AR38").Value = 67
$ws.Range("AX38").Value = 21
$ws.Range("AY38").Value = 34
$ws.Range("BB38").Value = 301
